$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round-half-to-even (banker's rounding) on the true binary value, matching
# the numpy/pandas semantics that produced the authoritative figures.
function RoundHalfEven($val, $digits) {
    $scale = [Math]::Pow(10, $digits)
    $scaled = $val * $scale
    $floor = [Math]::Floor($scaled)
    $diff = $scaled - $floor
    if ($diff -lt 0.5) {
        $rounded = $floor
    } elseif ($diff -gt 0.5) {
        $rounded = $floor + 1
    } else {
        if (($floor % 2) -eq 0) {
            $rounded = $floor
        } else {
            $rounded = $floor + 1
        }
    }
    return $rounded / $scale
}

# Data occupies rows 2..54 (row 1 is the header row).
$firstRow = 2
$lastRow = 54

# Column layout: A=index, B=month, C=week, D=monthly_mean, E=weekly_mean, F=weekly_share
$colMonth = 2
$colMonthlyMean = 4
$colWeeklyMean = 5
$colWeeklyShare = 6

# First pass: count how many weekly rows belong to each month.
$weeksPerMonth = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $month = $ws.Cells.Item($r, $colMonth).Value2
    if ($weeksPerMonth.ContainsKey($month)) {
        $weeksPerMonth[$month] = $weeksPerMonth[$month] + 1
    } else {
        $weeksPerMonth[$month] = 1
    }
}

# Second pass: monthly_mean is re-derived by spreading it across the weeks in
# that month (divide by week count), and weekly_share is recomputed from the
# (unchanged) weekly_mean divided by the new monthly_mean.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $month = $ws.Cells.Item($r, $colMonth).Value2
    $weeks = $weeksPerMonth[$month]

    $oldMonthlyMean = $ws.Cells.Item($r, $colMonthlyMean).Value2
    $weeklyMean = $ws.Cells.Item($r, $colWeeklyMean).Value2

    $newMonthlyMean = RoundHalfEven ($oldMonthlyMean / $weeks) 3
    $newWeeklyShare = RoundHalfEven ($weeklyMean / $newMonthlyMean) 3

    $ws.Cells.Item($r, $colMonthlyMean).Value = $newMonthlyMean
    $ws.Cells.Item($r, $colWeeklyShare).Value = $newWeeklyShare
}
